# Inserts a new data row at row 50 (shifting existing rows 50:183 down to 51:184)
# and populates it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 50; this shifts rows 50:183 -> 51:184
$ws.Rows("50:50").Insert()

# Populate the newly inserted row 50 with the new record
$ws.Range("A50").Value = 7
$ws.Range("B50").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C50").Value = "Ñuble"
$ws.Range("D50").Value = 44607
$ws.Range("E50").Value = 16
$ws.Range("F50").Value = 100112006
$ws.Range("G50").Value = "Repollo"
$ws.Range("H50").Value = "Crespo record"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 200
$ws.Range("K50").Value = 700
$ws.Range("L50").Value = 750
$ws.Range("M50").Value = 725
$ws.Range("N50").Value = "$/unidad"
$ws.Range("O50").Value = "Provincia de Diguillín"
$ws.Range("P50").Value = 725
$ws.Range("Q50").Value = 1
$ws.Range("R50").Value = "Hortaliza"
